# Swap the "Steps" / "Expected Results" content between the TC2 and TC4
# test-case blocks (TC3 stays the same), per commit "From v1.1 to v1.1.1".
#
# Layout on the sheet (ActiveSheet):
#   TC2 block: B18 = Steps, D18 = Expected Results
#   TC3 block: B25 = Steps, D25 = Expected Results   (unchanged)
#   TC4 block: B32 = Steps, D32 = Expected Results

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc2Steps   = $ws.Range("B18").Value2
$tc2Results = $ws.Range("D18").Value2
$tc4Steps   = $ws.Range("B32").Value2
$tc4Results = $ws.Range("D32").Value2

$ws.Range("B18").Value = $tc4Steps
$ws.Range("D18").Value = $tc4Results

$ws.Range("B32").Value = $tc2Steps
$ws.Range("D32").Value = $tc2Results
